$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(516).Delete()
